$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A3").Value = 2
$ws.Range("B3").Value = "veronica"
$ws.Range("C3").Value = "soloduha"
$ws.Range("D3").Value = "veronica"
$ws.Range("E3").Value = "soloduha"
$ws.Range("F3").Value = "Tue Apr 19 22:14:07 2022"
